$d = $word.ActiveDocument

# The second paragraph holds a Word field: { m:link Var.name }.
# Replace the field (fldChar begin/instrText runs/fldChar end) with
# plain literal-text runs spelling out "{m:link Var.name}" while
# preserving the bookmark and the colored formatting on "link"/"Var".
$p = $d.Paragraphs(2)

$xml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?>' +
  '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
  '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
  '<pkg:xmlData>' +
  '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
  '<w:body>' +
  '<w:p w:rsidP="00F5495F" w:rsidR="00C52979" w:rsidRDefault="00C52979">' +
  '<w:r><w:t>{</w:t></w:r>' +
  '<w:r><w:t>m</w:t></w:r>' +
  '<w:r><w:t>:</w:t></w:r>' +
  '<w:r><w:rPr><w:color w:themeColor="accent6" w:themeShade="BF" w:val="E36C0A"/></w:rPr><w:t>link</w:t></w:r>' +
  '<w:bookmarkStart w:id="0" w:name="_GoBack"/>' +
  '<w:bookmarkEnd w:id="0"/>' +
  '<w:r><w:rPr><w:color w:themeColor="accent6" w:themeShade="BF" w:val="E36C0A"/></w:rPr><w:t>Var</w:t></w:r>' +
  '<w:r><w:t xml:space="preserve">.name}</w:t></w:r>' +
  '</w:p>' +
  '</w:body>' +
  '</w:document>' +
  '</pkg:xmlData>' +
  '</pkg:part>' +
  '</pkg:package>'

$p.Range.InsertXML($xml)
